$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'296.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-2.16%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'31.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-1.76%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.115"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.79%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07372"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.33%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.703"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-1.14%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.736"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.17%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.631"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'11.05%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9172"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.28%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1671"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.58%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07162"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-4.41%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07944"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.04%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'0.46%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09910"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.01%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001489"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.80%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006098"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-6.11%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'-0.47%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'0.10%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3273"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.84%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1332"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.40%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.547"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'6.23%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04618"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.87%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.1549"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-4.88%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-0.34%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004428"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.23%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001298"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.32%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001873"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'7.41%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01685"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'1.91%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04395"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.75%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007216"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.39%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1328"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.93%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002136"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-8.45%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'-16.98%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006021"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'1.928"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'1.87%"
$ws.Range("E46").Style = "Normal"
